# Add a new "COST PRICE" column to the Sheet1 import header row, just
# after the existing "MU" column (P1), matching the commit message:
# "Added CostPrice entry in create, export, import, undo, update flows".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing header cell (P1, style "MU")
# onto the new header cell Q1, so the new header looks consistent with
# the rest of the header row (bold, centered, yellow fill, bordered).
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text / shared string value.
$ws.Range("Q1").Value = "COST PRICE"

# Reflect the author's final cell selection after making the edit.
$null = $ws.Range("R2").Select()
